$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove two rows that were dropped from the sample set ---
# Row 28 = "SC 92", row 26 = "RM 232" in the original layout.
# Delete the lower one first so the other row index stays valid.
$ws.Rows(28).Delete()
$ws.Rows(26).Delete()

# --- Toggle individual "missing value" cells (post-deletion row numbers) ---

# F5 (RM 14): value removed -> becomes missing
$ws.Range("F5").Value = "'"

# D6 (RM 21): value restored
$ws.Range("D6").Value = -14.2

# D8 (RM 38): value removed -> becomes missing
$ws.Range("D8").Value = "'"

# F11 (RM 58): value restored
$ws.Range("F11").Value = 17.65

# D19 (RM 125): value restored
$ws.Range("D19").Value = -15.5
# F19 (RM 125): value removed -> becomes missing
$ws.Range("F19").Value = "'"

# D21 (RM 135): value removed -> becomes missing
$ws.Range("D21").Value = "'"

# D23 (RM 140): value restored
$ws.Range("D23").Value = -13.9
# F23 (RM 140): value restored
$ws.Range("F23").Value = 16.48

# F25 (RM 145): value restored
$ws.Range("F25").Value = 16.6

# C26 (SC 5): value removed -> becomes missing
$ws.Range("C26").Value = "'"

# C27 (SC 101): value restored
$ws.Range("C27").Value = 10
# D27 (SC 101): value removed -> becomes missing
$ws.Range("D27").Value = "'"
# F27 (SC 101): value removed -> becomes missing
$ws.Range("F27").Value = "'"

# C29 (SC 119): value removed -> becomes missing
$ws.Range("C29").Value = "'"
# D29 (SC 119): value restored
$ws.Range("D29").Value = -13
# F29 (SC 119): value removed -> becomes missing
$ws.Range("F29").Value = "'"

# F30 (SC 120): value restored
$ws.Range("F30").Value = 16.89

# F33 (SC 232): value restored
$ws.Range("F33").Value = 17.53
